# BOM.xlsx edit: add Silica Gel line item, add PCB price/hyperlink, add Total Cost row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 13 (the "Acrylic Case:" section) to host the
# new "Silica Gel" line item. This shifts every row from 13 downward by one.
$ws.Rows.Item(12).Insert()

# New row 12: Silica Gel, 1kg, qty 1, eBay, url, 8.95
$ws.Range("B12").Value = "Silica Gel"
$ws.Range("C12").Value = "1kg"
$ws.Range("D12").Value = 1
$ws.Range("F12").Value = "eBay"
$ws.Range("G12").Value = "http://www.ebay.co.uk/itm/1000g-1kg-BONUS-100g-BAG-Silica-Gel-Desiccant-Self-Indicating-Loose-/111737656127?hash=item1a0415533f:g:340AAOSwu4BVm9lH"
$ws.Range("H12").Value = 8.95

# Row 19 (was row 18 before the insert) is the "PCB / OSH Park" line; add
# the quantity and price for it.
$ws.Range("E19").Value = 1
$ws.Range("H19").Value = 15

# New totals row at the bottom of the sheet (written before the "USD 20 -
# for 3." note below so the shared-string table fills in the same order
# as the authored workbook).
$ws.Range("A45").Value = "Total Cost"
$ws.Range("H45").Formula = "=SUM(H4:H44)"

$ws.Range("K19").Value = "USD 20 - for 3."

# Turn the OSH Park URL already in G19 into a real hyperlink (same target
# text, now styled + linked).
$ws.Hyperlinks.Add($ws.Range("G19"), "https://oshpark.com/shared_projects/nOVDuNCE")
$ws.Range("G19").Style = "Hyperlink"

# Restore the selection the author ended up with.
[void]$ws.Range("K20").Select()
